# Apply odds updates to the "Jogos da Semana" FlashScore workbook and
# remove the last match row (row 8) that is no longer part of the report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 updates ---
$ws.Range("K4").Value2  = 2.25
$ws.Range("AC4").Value2 = 11
$ws.Range("AD4").Value2 = 7
$ws.Range("AI4").Value2 = 17
$ws.Range("AQ4").Value2 = 51
$ws.Range("AW4").Value2 = 34

# --- Row 5 updates ---
$ws.Range("H5").Value2  = 3
$ws.Range("I5").Value2  = 3.1
$ws.Range("L5").Value2  = 4
$ws.Range("O5").Value2  = 1.53
$ws.Range("P5").Value2  = 2.5
$ws.Range("U5").Value2  = 2.1
$ws.Range("V5").Value2  = 1.67
$ws.Range("Y5").Value2  = 11
$ws.Range("Z5").Value2  = 26
$ws.Range("AT5").Value2 = 81
$ws.Range("AZ5").Value2 = 351

# --- Row 7 updates ---
$ws.Range("G7").Value2  = 3.2
$ws.Range("I7").Value2  = 2.5
$ws.Range("M7").Value2  = 1.11
$ws.Range("N7").Value2  = 6.5
$ws.Range("X7").Value2  = 13
$ws.Range("AH7").Value2 = 11
$ws.Range("AK7").Value2 = 26
$ws.Range("AM7").Value2 = 4.75
$ws.Range("AY7").Value2 = 101

# --- Remove the last data row (row 8, Uruguay Progreso - Penarol) ---
$ws.Rows.Item(8).Delete()
